$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, copying the formatting of the neighboring
# header cell (G1) so the new column matches the existing header style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the new data values for the Save column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
